# Generate Report for Handoff
# Replace the old localization-run identifiers / hashes / timestamps with the
# new ones produced by this handoff run (new source-file GUID + new xliff
# content hashes + refreshed "ready for handoff" timestamps).

$wb = $excel.ActiveWorkbook

$newGuid = "4654bfce-3af4-4efb-b4b3-6538e811e7ce"
$newZhHash = "7c31b78e56a297e381cf7153862ef697fa4b1cda"
$newDeHash = "7c31b78e56a297e381cf7153862ef697fa4b1cda"

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
foreach ($h in $wsOverview.Hyperlinks) {
    $h.TextToDisplay = "e2e\$newGuid.md"
}
$wsOverview.Range("G2").Value = "2016-09-04 13:02:11"

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = "$newGuid.md"
foreach ($h in $wsZh.Hyperlinks) {
    $h.TextToDisplay = "$newGuid.md"
}
$wsZh.Range("G2").Value = "$newGuid.$newZhHash.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-09-04 13:02:01"

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = "$newGuid.md"
foreach ($h in $wsDe.Hyperlinks) {
    $h.TextToDisplay = "$newGuid.md"
}
$wsDe.Range("G2").Value = "$newGuid.$newDeHash.de-de.xlf"
$wsDe.Range("H2").Value = "2016-09-04 13:02:11"
